$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 325
$ws.Range("I18").Value = 325
$ws.Range("K18").Value = 325
$ws.Range("M18").Value = -41

$ws.Range("H33").Value = 3875759
$ws.Range("J33").Value = 395
$ws.Range("L33").Value = 395
$ws.Range("N33").Value = -853

$ws.Range("H116").Value = 6610.7407
$ws.Range("I116").Value = 5788.1816
$ws.Range("J116").Value = 7176.25
$ws.Range("K116").Value = 5788.1816
$ws.Range("L116").Value = 7176.25
$ws.Range("M116").Value = -2346.1816
$ws.Range("N116").Value = -14060.25

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H137").Value = 3744.59
$ws.Range("I137").Value = 3828.7693
$ws.Range("K137").Value = 11486.3079
$ws.Range("M137").Value = -8936.3079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 180
$ws.Range("I5").Value = 193.125
$ws.Range("J5").Value = 75
$ws.Range("K5").Value = 193.125
$ws.Range("L5").Value = 75
$ws.Range("M5").Value = -81.125
$ws.Range("N5").Value = -299

$ws.Range("H61").Value = 6649.8335
$ws.Range("I61").Value = 3530.111
$ws.Range("K61").Value = 3530.111
$ws.Range("M61").Value = -3318.111

$ws.Range("H74").Value = 4566.125
$ws.Range("I74").Value = 1596.8889
$ws.Range("J74").Value = 8383.714
$ws.Range("K74").Value = 1596.8889
$ws.Range("L74").Value = 8383.714
$ws.Range("M74").Value = -722.8888999999999
$ws.Range("N74").Value = -10131.714

$ws.Range("H77").Value = 4566.125
$ws.Range("I77").Value = 1596.8889
$ws.Range("J77").Value = 8383.714
$ws.Range("K77").Value = 7984.4445
$ws.Range("L77").Value = 41918.57
$ws.Range("M77").Value = -3616.4445
$ws.Range("N77").Value = -50654.57

$ws.Range("H110").Value = 3698.3333
$ws.Range("J110").Value = 4875
$ws.Range("L110").Value = 4875
$ws.Range("N110").Value = -8965

$ws.Range("H122").Value = 2131.5938
$ws.Range("I122").Value = 2170.5386
$ws.Range("K122").Value = 6511.6158
$ws.Range("M122").Value = -4061.6158

$ws.Range("H136").Value = 6649.8335
$ws.Range("I136").Value = 3530.111
$ws.Range("K136").Value = 10590.333
$ws.Range("M136").Value = -8040.332999999999

$ws.Range("H139").Value = 96999
$ws.Range("J139").Value = 96999
$ws.Range("L139").Value = 96999
$ws.Range("N139").Value = -107279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 180
$ws.Range("I4").Value = 193.125
$ws.Range("J4").Value = 75
$ws.Range("K4").Value = 193.125
$ws.Range("L4").Value = 75
$ws.Range("M4").Value = -78.125
$ws.Range("N4").Value = -305

$ws.Range("H11").Value = 478.85715
$ws.Range("J11").Value = 2944.5
$ws.Range("L11").Value = 2944.5
$ws.Range("N11").Value = -3224.5

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H105").Value = 2755.1904
$ws.Range("I105").Value = 3005.2307
$ws.Range("K105").Value = 3005.2307
$ws.Range("M105").Value = -1258.2307

$ws.Range("H134").Value = 7761.3076
$ws.Range("I134").Value = 6785.591
$ws.Range("K134").Value = 20356.773
$ws.Range("M134").Value = -17821.773

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2490.65
$ws.Range("I16").Value = 2515.8125
$ws.Range("J16").Value = 2390
$ws.Range("K16").Value = 2515.8125
$ws.Range("L16").Value = 2390
$ws.Range("M16").Value = -2228.8125
$ws.Range("N16").Value = -2964

$ws.Range("H22").Value = 359.77274
$ws.Range("I22").Value = 318
$ws.Range("J22").Value = 547.75
$ws.Range("K22").Value = 318
$ws.Range("L22").Value = 547.75
$ws.Range("M22").Value = 32
$ws.Range("N22").Value = -1247.75

$ws.Range("H31").Value = 1764.1482
$ws.Range("I31").Value = 1649
$ws.Range("K31").Value = 1649
$ws.Range("M31").Value = -1354

$ws.Range("H34").Value = 1764.1482
$ws.Range("I34").Value = 1649
$ws.Range("K34").Value = 1649
$ws.Range("M34").Value = -1447

$ws.Range("H107").Value = 744.25
$ws.Range("I107").Value = 542.4545000000001
$ws.Range("K107").Value = 542.4545000000001
$ws.Range("M107").Value = 1377.5455

$ws.Range("H113").Value = 2490.65
$ws.Range("I113").Value = 2515.8125
$ws.Range("J113").Value = 2390
$ws.Range("K113").Value = 2515.8125
$ws.Range("L113").Value = 2390
$ws.Range("M113").Value = -345.8125
$ws.Range("N113").Value = -6730

$ws.Range("H134").Value = 7809.5
$ws.Range("I134").Value = 5899.027
$ws.Range("K134").Value = 17697.081
$ws.Range("M134").Value = -15162.081

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 33335552
$ws.Range("I132").Value = 3323
$ws.Range("K132").Value = 29907
$ws.Range("M132").Value = -27377

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 119.5
$ws.Range("J2").Value = 191.83333
$ws.Range("L2").Value = 191.83333
$ws.Range("N2").Value = -417.83333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1783.1666
$ws.Range("I7").Value = 1633
$ws.Range("J7").Value = 1933.3334
$ws.Range("K7").Value = 1633
$ws.Range("L7").Value = 1933.3334
$ws.Range("M7").Value = -1521
$ws.Range("N7").Value = -2157.3334

$ws.Range("H22").Value = 5779
$ws.Range("J22").Value = 8332
$ws.Range("L22").Value = 8332
$ws.Range("N22").Value = -8922

$ws.Range("H27").Value = 5779
$ws.Range("J27").Value = 8332
$ws.Range("L27").Value = 8332
$ws.Range("N27").Value = -8546

$ws.Range("H36").Value = 78899.60000000001
$ws.Range("J36").Value = 78899.60000000001
$ws.Range("L36").Value = 78899.60000000001
$ws.Range("N36").Value = -80023.60000000001

$ws.Range("H100").Value = 587125.8
$ws.Range("I100").Value = 2436.5
$ws.Range("J100").Value = 1012354.4
$ws.Range("K100").Value = 2436.5
$ws.Range("L100").Value = 1012354.4
$ws.Range("M100").Value = -1895.5
$ws.Range("N100").Value = -1013436.4

$ws.Range("H126").Value = 1783.1666
$ws.Range("I126").Value = 1633
$ws.Range("J126").Value = 1933.3334
$ws.Range("K126").Value = 4899
$ws.Range("L126").Value = 5800.0002
$ws.Range("M126").Value = -2429
$ws.Range("N126").Value = -10740.0002

$ws.Range("H132").Value = 1715556
$ws.Range("J132").Value = 6164.4443
$ws.Range("L132").Value = 18493.3329
$ws.Range("N132").Value = -23553.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 16239.516
$ws.Range("I132").Value = 10959.358
$ws.Range("K132").Value = 32878.074
$ws.Range("M132").Value = -30348.074

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H141").Value = 80707
$ws.Range("J141").Value = 80707
$ws.Range("L141").Value = 80707
$ws.Range("N141").Value = -91067
